# Third Iteration - Testing
# The document ends with a trailing empty paragraph (right after the
# "After you make the diagrams..." paragraph, before the section break).
# That empty paragraph must be removed so the document's last paragraph
# becomes the one containing the "After you make the diagrams..." text.

$d = $word.ActiveDocument

# The very last paragraph in the document is the empty one we need to drop.
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)

# Deleting its Range removes the paragraph mark along with it, collapsing
# it out of the document entirely (rather than just clearing its text).
$lastParagraph.Range.Delete()
